# Update column F (dSF) values for rows 2-15 per repull/recalculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    3  = 5
    4  = -4
    5  = 0
    6  = -2
    7  = -5
    9  = -3
    10 = -1
    11 = 3
    12 = 1
    13 = 1
    14 = -3
    15 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
